$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 500047
$ws.Range("I9").Value = 500047
$ws.Range("K9").Value = 500047
$ws.Range("M9").Value = -499878
$ws.Range("H15").Value = 1839.7059
$ws.Range("I15").Value = 1839.7059
$ws.Range("K15").Value = 5519.1177
$ws.Range("M15").Value = -5350.1177
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1000
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -1350
$ws.Range("H92").Value = 1692.4615
$ws.Range("I92").Value = 1739.2222
$ws.Range("J92").Value = 1587.25
$ws.Range("K92").Value = 1739.2222
$ws.Range("L92").Value = 1587.25
$ws.Range("M92").Value = -491.2221999999999
$ws.Range("N92").Value = -4083.25
$ws.Range("H135").Value = 1265.8
$ws.Range("I135").Value = 999.1429000000001
$ws.Range("K135").Value = 8992.286100000001
$ws.Range("M135").Value = -6457.286100000001
$ws.Range("H137").Value = 1963.6111
$ws.Range("I137").Value = 1856.4
$ws.Range("K137").Value = 5569.200000000001
$ws.Range("M137").Value = -3019.200000000001
$ws.Range("H138").Value = 9061.125
$ws.Range("J138").Value = 9331.866
$ws.Range("L138").Value = 27995.598
$ws.Range("N138").Value = -38275.598

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11030.556
$ws.Range("I32").Value = 8457.210999999999
$ws.Range("J32").Value = 25000.143
$ws.Range("K32").Value = 8457.210999999999
$ws.Range("L32").Value = 25000.143
$ws.Range("M32").Value = -8170.210999999999
$ws.Range("N32").Value = -25574.143
$ws.Range("H45").Value = 2998
$ws.Range("I45").Value = 2998
$ws.Range("K45").Value = 2998
$ws.Range("M45").Value = -2621
$ws.Range("H46").Value = 13096.667
$ws.Range("I46").Value = 10069
$ws.Range("J46").Value = 19152
$ws.Range("K46").Value = 10069
$ws.Range("L46").Value = 19152
$ws.Range("M46").Value = -9750
$ws.Range("N46").Value = -19790
$ws.Range("H61").Value = 3849.75
$ws.Range("I61").Value = 3849.75
$ws.Range("K61").Value = 3849.75
$ws.Range("M61").Value = -3637.75
$ws.Range("H74").Value = 8753
$ws.Range("I74").Value = 11004
$ws.Range("K74").Value = 11004
$ws.Range("M74").Value = -10130
$ws.Range("H77").Value = 8753
$ws.Range("I77").Value = 11004
$ws.Range("K77").Value = 55020
$ws.Range("M77").Value = -50652
$ws.Range("H122").Value = 2785.8823
$ws.Range("J122").Value = 3013.5
$ws.Range("L122").Value = 9040.5
$ws.Range("N122").Value = -13940.5
$ws.Range("H136").Value = 3849.75
$ws.Range("I136").Value = 3849.75
$ws.Range("K136").Value = 11549.25
$ws.Range("M136").Value = -8999.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 80000
$ws.Range("J74").Value = 80000
$ws.Range("L74").Value = 80000
$ws.Range("N74").Value = -81872
$ws.Range("H77").Value = 80000
$ws.Range("J77").Value = 80000
$ws.Range("L77").Value = 240000
$ws.Range("N77").Value = -249360
$ws.Range("H86").Value = 10700
$ws.Range("I86").Value = 5500
$ws.Range("K86").Value = 5500
$ws.Range("M86").Value = -4377
$ws.Range("H89").Value = 10700
$ws.Range("I89").Value = 5500
$ws.Range("K89").Value = 27500
$ws.Range("M89").Value = -21884
$ws.Range("H94").Value = 5341.5
$ws.Range("I94").Value = 5957
$ws.Range("K94").Value = 5957
$ws.Range("M94").Value = -5506
$ws.Range("H134").Value = 1517.8125
$ws.Range("I134").Value = 1345
$ws.Range("J134").Value = 2266.6667
$ws.Range("K134").Value = 4035
$ws.Range("L134").Value = 6800.000100000001
$ws.Range("M134").Value = -1500
$ws.Range("N134").Value = -11870.0001
$ws.Range("H141").Value = 199999
$ws.Range("J141").Value = 199999
$ws.Range("L141").Value = 199999
$ws.Range("N141").Value = -210359

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 1700
$ws.Range("I32").Value = 1700
$ws.Range("K32").Value = 1700
$ws.Range("M32").Value = -1384
$ws.Range("H86").Value = 7540.857
$ws.Range("I86").Value = 6878.2
$ws.Range("K86").Value = 6878.2
$ws.Range("M86").Value = -5755.2
$ws.Range("H89").Value = 7540.857
$ws.Range("I89").Value = 6878.2
$ws.Range("K89").Value = 34391
$ws.Range("M89").Value = -28775
$ws.Range("H105").Value = 1574.5
$ws.Range("I105").Value = 1574.5
$ws.Range("K105").Value = 1574.5
$ws.Range("M105").Value = 172.5
$ws.Range("H122").Value = 2476.6667
$ws.Range("I122").Value = 2441.125
$ws.Range("K122").Value = 7323.375
$ws.Range("M122").Value = -4873.375
$ws.Range("H134").Value = 2198.9565
$ws.Range("I134").Value = 1705.6111
$ws.Range("K134").Value = 5116.8333
$ws.Range("M134").Value = -2581.8333
$ws.Range("H141").Value = 699333
$ws.Range("J141").Value = 699333
$ws.Range("L141").Value = 699333
$ws.Range("N141").Value = -709693

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 14349.75
$ws.Range("J55").Value = 14349.75
$ws.Range("L55").Value = 43049.25
$ws.Range("N55").Value = -43403.25
$ws.Range("H118").Value = 4904.7144
$ws.Range("I118").Value = 3999.5
$ws.Range("K118").Value = 11998.5
$ws.Range("M118").Value = -10755.5
$ws.Range("H141").Value = 5469.5
$ws.Range("I141").Value = 5712
$ws.Range("K141").Value = 17136
$ws.Range("M141").Value = -11956

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 10000
$ws.Range("M70").Value = -9730
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 10000
$ws.Range("M73").Value = -9064
$ws.Range("H132").Value = 2881.6
$ws.Range("I132").Value = 2020.3334
$ws.Range("K132").Value = 6061.0002
$ws.Range("M132").Value = -3531.0002

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2892
$ws.Range("I68").Value = 2892
$ws.Range("K68").Value = 2892
$ws.Range("M68").Value = -2143
$ws.Range("H71").Value = 2892
$ws.Range("I71").Value = 2892
$ws.Range("K71").Value = 14460
$ws.Range("M71").Value = -10716

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2875.6
$ws.Range("I81").Value = 2195.111
$ws.Range("K81").Value = 4390.222
$ws.Range("M81").Value = -3329.222
$ws.Range("H84").Value = 2875.6
$ws.Range("I84").Value = 2195.111
$ws.Range("K84").Value = 21951.11
$ws.Range("M84").Value = -16647.11
$ws.Range("H103").Value = 350000
$ws.Range("J103").Value = 350000
$ws.Range("L103").Value = 350000
$ws.Range("N103").Value = -352344
$ws.Range("H107").Value = 684.9231
$ws.Range("I107").Value = 564.1429000000001
$ws.Range("J107").Value = 825.8333
$ws.Range("K107").Value = 1692.4287
$ws.Range("L107").Value = 2477.4999
$ws.Range("M107").Value = 227.5712999999998
$ws.Range("N107").Value = -6317.4999
